# Fruta / hortaliza, semanal
# The weekly refresh re-shuffles the per-record fields (Fecha, Volumen,
# Precio minimo/maximo/promedio ponderado, Origen, Precio $/Kg) across the
# existing data rows (2-17) of the active sheet. Columns A, B, C, E-L, Q and T
# are unaffected because their values are already identical on every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: target row -> source row (values are read from the "before" state
# of the source row and written into the target row).
$rowMap = @{
    2  = 5
    3  = 4
    4  = 6
    5  = 17
    6  = 9
    7  = 7
    8  = 2
    9  = 15
    10 = 13
    11 = 12
    12 = 11
    13 = 3
    14 = 14
    15 = 10
    16 = 8
    17 = 16
}

# Snapshot the "before" values for the columns that move, keyed by row number,
# so that reads always reflect the original workbook state even after we
# start writing.
$cols = @("D", "M", "N", "O", "P", "R", "S")
$snapshot = @{}
foreach ($r in 2..17) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($r in 2..17) {
    $srcRow = $rowMap[$r]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value2 = $srcVals[$col]
    }
}
